$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "20 Best jobs in Colorado Springs, CO (Hiring Now!) | SimplyHired"
$ws.Range("B2").Value = "f9dca67454f647c9ac710a500a0bbebb@sentry.indeed.com`n"
$ws.Range("C2").Value = 629
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"

$ws.Range("A3").Value = "All Jobs in Colorado Springs, CO - Apply Now | CareerBuilder"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = 430
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 18
$ws.Range("F3").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"

$ws.Range("A4").Value = "`$60,000 Jobs, Employment in Colorado Springs, CO | Indeed.com"
$ws.Range("B4").Value = "0252655a41544fd28ae41f8b8ff36917@sentry.indeed.com`n"
$ws.Range("C4").Value = 1471
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 22
$ws.Range("F4").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"

$ws.Range("A5").Value = "Jobs, Employment in Colorado Springs, CO | Indeed.com"
$ws.Range("B5").Value = "0252655a41544fd28ae41f8b8ff36917@sentry.indeed.com`n"
$ws.Range("C5").Value = 1553
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 18
$ws.Range("F5").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"

$ws.Range("A6").Value = "Apache Tomcat - Error report"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"

$ws.Range("A7").Value = "`$34k-`$81k Jobs in Colorado Springs, CO | ZipRecruiter"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = 1452
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 12
$ws.Range("F7").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"

$ws.Range("A8").Value = "LinkedIn Job Search: Find US Jobs, Internships, Jobs Near Me"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = 475
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"

$ws.Range("A9").Value = "Find a Job | Careers in Colorado Springs, CO | Monster"
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = 440
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"

$ws.Range("A10").Value = "City Jobs and Careers | Colorado Springs"
$ws.Range("B10").Value = "nr@context`nnr@id`nnr@original`nhr@springsgov.com`nwebmaster@springsgov.com`n"
$ws.Range("C10").Value = 59
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"

$ws.Range("A11").Value = "Jobs in Colorado Springs, Co Now Hiring | Snagajob"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = 431
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = "https://www.snagajob.com/search/w-colorado+springs,+co"
